# Apply CTA violent crime YTD update for 2023-12-16
# Updates citywide totals, per-neighborhood 'By Neighborhood' summary sheet,
# and individual neighborhood detail sheets with revised year-to-date crime counts.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("D2").Value = 92
$ws.Range("H2").Value = 108
$ws.Range("B3").Value = 77
$ws.Range("J3").Value = 237
$ws.Range("C6").Value = 489
$ws.Range("D6").Value = 421
$ws.Range("E6").Value = 487
$ws.Range("F6").Value = 553
$ws.Range("G6").Value = 439
$ws.Range("I6").Value = 506
$ws.Range("B7").Value = 517
$ws.Range("C7").Value = 646
$ws.Range("D7").Value = 660
$ws.Range("E7").Value = 720
$ws.Range("F7").Value = 801
$ws.Range("G7").Value = 673
$ws.Range("H7").Value = 727
$ws.Range("I7").Value = 844
$ws.Range("J7").Value = 808

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("C6").Value = 35
$ws.Range("C7").Value = 40

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("C4").Value = 8
$ws.Range("C5").Value = 10

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("G5").Value = 14
$ws.Range("G6").Value = 22

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("H2").Value = 11
$ws.Range("B3").Value = 4
$ws.Range("D6").Value = 24
$ws.Range("F6").Value = 38
$ws.Range("B7").Value = 36
$ws.Range("D7").Value = 48
$ws.Range("F7").Value = 58
$ws.Range("H7").Value = 46

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I7").Value = 10
$ws.Range("B28").Value = 36
$ws.Range("D28").Value = 48
$ws.Range("F28").Value = 58
$ws.Range("H28").Value = 46
$ws.Range("J35").Value = 8
$ws.Range("C36").Value = 40
$ws.Range("C41").Value = 10
$ws.Range("F47").Value = 18
$ws.Range("D53").Value = 73
$ws.Range("I68").Value = 6
$ws.Range("E70").Value = 19
$ws.Range("E74").Value = 7
$ws.Range("I74").Value = 20
$ws.Range("D76").Value = 15
$ws.Range("C77").Value = 25
$ws.Range("G86").Value = 22
$ws.Range("E95").Value = 6
$ws.Range("I95").Value = 6
$ws.Range("B98").Value = 517
$ws.Range("C98").Value = 646
$ws.Range("D98").Value = 660
$ws.Range("E98").Value = 720
$ws.Range("F98").Value = 801
$ws.Range("G98").Value = 673
$ws.Range("H98").Value = 727
$ws.Range("I98").Value = 844
$ws.Range("J98").Value = 808

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("D6").Value = 44
$ws.Range("D7").Value = 73

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J3").Value = 2
$ws.Range("J6").Value = 8

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("D2").Value = 3
$ws.Range("D7").Value = 15

$ws = $wb.Worksheets.Item('River North')
$ws.Range("E5").Value = 5
$ws.Range("I5").Value = 13
$ws.Range("E6").Value = 7
$ws.Range("I6").Value = 20

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("C6").Value = 16
$ws.Range("C7").Value = 25

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I5").Value = 5
$ws.Range("I6").Value = 10

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("F5").Value = 11
$ws.Range("F6").Value = 18

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("E4").Value = 17
$ws.Range("E5").Value = 19

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("H5").Value = 2
$ws.Range("H6").Value = 6

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("D4").Value = 3
$ws.Range("G4").Value = 5
$ws.Range("D5").Value = 6
$ws.Range("G5").Value = 6
